# Update the "Scaled" sheet (numeric, scaled 0-100 columns D:J)
$wb = $excel.ActiveWorkbook

$scaled = $wb.Worksheets.Item("Scaled")

$scaled.Range("D2").Value = 84.0
$scaled.Range("E2").Value = 70.0
$scaled.Range("F2").Value = 73.0
$scaled.Range("G2").Value = 0.0
$scaled.Range("H2").Value = 96.0
$scaled.Range("I2").Value = 90.0
$scaled.Range("J2").Value = 413.0

$scaled.Range("D3").Value = 99.0
$scaled.Range("E3").Value = 100.0
$scaled.Range("F3").Value = 100.0
$scaled.Range("G3").Value = 37.0
$scaled.Range("H3").Value = 98.0
$scaled.Range("I3").Value = 90.0
$scaled.Range("J3").Value = 524.0

$scaled.Range("D4").Value = 96.0
$scaled.Range("E4").Value = 87.0
$scaled.Range("F4").Value = 98.0
$scaled.Range("G4").Value = 100.0
$scaled.Range("H4").Value = 83.0
$scaled.Range("I4").Value = 5.0
$scaled.Range("J4").Value = 469.0

$scaled.Range("D5").Value = 85.0
$scaled.Range("E5").Value = 66.0
$scaled.Range("F5").Value = 63.0
$scaled.Range("G5").Value = 0.0
$scaled.Range("H5").Value = 65.0
$scaled.Range("I5").Value = 82.0
$scaled.Range("J5").Value = 361.0

$scaled.Range("D6").Value = 71.0
$scaled.Range("E6").Value = 100.0
$scaled.Range("F6").Value = 68.0
$scaled.Range("G6").Value = 55.0
$scaled.Range("H6").Value = 85.0
$scaled.Range("I6").Value = 100.0
$scaled.Range("J6").Value = 479.0

# Update the "Raw" sheet (text columns D:I — MDL/SPT/HRP are numeric-looking
# text, so force Text format first so Excel doesn't silently coerce them to
# numbers; SDC/PLK/2MR are clock-style strings that stay text on their own).
$raw = $wb.Worksheets.Item("Raw")

$raw.Range("D2:F6").NumberFormat = "@"

$raw.Range("D2").Value = "262"
$raw.Range("E2").Value = "7.0"
$raw.Range("F2").Value = "262"
$raw.Range("G2").Value = "4:47"
$raw.Range("H2").Value = "3:24"
$raw.Range("I2").Value = "15:34"

$raw.Range("D3").Value = "338"
$raw.Range("E3").Value = "10.0"
$raw.Range("F3").Value = "338"
$raw.Range("G3").Value = "2:51"
$raw.Range("H3").Value = "3:34"
$raw.Range("I3").Value = "15:24"

$raw.Range("D4").Value = "324"
$raw.Range("E4").Value = "8.7"
$raw.Range("F4").Value = "324"
$raw.Range("G4").Value = "1:29"
$raw.Range("H4").Value = "2:47"
$raw.Range("I4").Value = "23:49"

$raw.Range("D5").Value = "276"
$raw.Range("E5").Value = "6.6"
$raw.Range("F5").Value = "276"
$raw.Range("G5").Value = "4:50"
$raw.Range("H5").Value = "1:38"
$raw.Range("I5").Value = "16:50"

$raw.Range("D6").Value = "147"
$raw.Range("E6").Value = "10.0"
$raw.Range("F6").Value = "147"
$raw.Range("G6").Value = "3:20"
$raw.Range("H6").Value = "2:48"
$raw.Range("I6").Value = "13:59"
